$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Trim the old sheet from 8 data rows down to 4 data rows (rows 6-8 no longer used) ---
$ws.Range("A6:C8").Delete() | Out-Null

# --- Row 1: header row (columns D-H are new; A-C get new header text) ---
$ws.Range("A1").Value = "Core Attribute"
$ws.Range("B1").Value = "http://dbpedia.org/ontology/deathPlace"
$ws.Range("C1").Value = "http://dbpedia.org/ontology/parent"
$ws.Range("D1").Value = "http://dbpedia.org/ontology/deathDate"
$ws.Range("E1").Value = "http://dbpedia.org/ontology/birthDate"
$ws.Range("F1").Value = "http://dbpedia.org/ontology/birthPlace"
$ws.Range("G1").Value = "http://dbpedia.org/ontology/associateStar"
$ws.Range("H1").Value = "http://dbpedia.org/ontology/fastestLap"

# Give the new header cells (D1:H1) the same bold/bordered style already used by A1:C1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:H1").PasteSpecial(-4122) | Out-Null

# --- Row 2 ---
$ws.Range("A2").Value = "http://dbpedia.org/resource/Giovanni_Francesco_Guidi_di_Bagno"
$ws.Range("B2").Value = "http://dbpedia.org/resource/Rome"
$ws.Range("C2").Value = "http://dbpedia.org/resource/Colonna_family"
$ws.Range("D2").Value = "http://dbpedia.org/resource/1641"
$ws.Range("E2").Value = "1578-10-04 "
$ws.Range("F2").Value = "http://dbpedia.org/resource/Florence http://dbpedia.org/resource/Grand_Duchy_of_Tuscany "
$ws.Range("G2").Value = "nan"
$ws.Range("H2").Value = "nan"

# --- Row 3 ---
$ws.Range("A3").Value = "http://dbpedia.org/resource/Giovanni_Doria"
$ws.Range("B3").Value = "http://dbpedia.org/resource/Palermo"
$ws.Range("C3").Value = "http://dbpedia.org/resource/Giovanni_Andrea_Doria"
$ws.Range("D3").Value = "http://dbpedia.org/resource/1642"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "nan"
$ws.Range("H3").Value = "nan"

# --- Row 4 ---
$ws.Range("A4").Value = "http://dbpedia.org/resource/Dick_Sheppard_(priest)"
$ws.Range("B4").Value = "http://dbpedia.org/resource/City_of_London "
$ws.Range("C4").Value = "http://dbpedia.org/resource/Edgar_Sheppard "
$ws.Range("D4").Value = "1937-10-31 "
$ws.Range("E4").Value = "http://dbpedia.org/resource/1880"
$ws.Range("F4").Value = "http://dbpedia.org/resource/Windsor"
$ws.Range("G4").Value = "nan"
$ws.Range("H4").Value = "nan"

# --- Row 5 ---
$ws.Range("A5").Value = "http://dbpedia.org/resource/Claus_Westermann"
$ws.Range("B5").Value = "http://dbpedia.org/resource/Heidelberg http://dbpedia.org/resource/Germany "
$ws.Range("C5").Value = "http://dbpedia.org/resource/Mother http://dbpedia.org/resource/Father "
$ws.Range("D5").Value = "2000-06-11 "
$ws.Range("E5").Value = "http://dbpedia.org/resource/1909"
$ws.Range("F5").Value = "http://dbpedia.org/resource/Berlin"
$ws.Range("G5").Value = "nan"
$ws.Range("H5").Value = "nan"
